$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing part values (U4,U3 op-amp swapped from TL072/C6961 to LM358DR2G/C7950) ---
$ws.Range("A8").Value = "LM358DR2G"
$ws.Range("D8").Value = "C7950"

# --- Normalize formatting: cells that used the old "default font" style (s=3) now use the
#     Arial/theme-color style (s=2) instead. Copy the format from a cell that already carries
#     that style so the underlying font/style table collapses the same way. ---
$ws.Range("A2").Copy()
$fixupCells = "B2,A4,B4,C4,D4,B5,B10"
foreach ($addr in $fixupCells.Split(",")) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# --- Append two new BOM rows ---
$ws.Range("A21").Value = "220pF"
$ws.Range("B21").Value = "C15"
$ws.Range("C21").Value = "C_0805_2012Metric"
$ws.Range("D21").Value = "C53172"

$ws.Range("A22").Value = "22pF"
$ws.Range("B22").Value = "C16"
$ws.Range("C22").Value = "C_0805_2012Metric"
$ws.Range("D22").Value = "C1804"

# match the formatting used throughout the rest of the BOM table
$ws.Range("A2").Copy()
$ws.Range("A21:D22").PasteSpecial(-4122)
$excel.CutCopyMode = 0
